$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new columns before column B (shifts existing B:V data to K:AE)
$ws.Columns("B:J").Insert()

# Fill in the new header dates (row 1), most-recent-first, continuing the
# existing weekly cadence after Jun_09
$ws.Range("B1").Value = "Sep_08"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("J1").Value = "Jun_16"

# Fill the new columns for the analyst rows with the "UN" (unrated / no
# action that period) placeholder, matching the existing data rows.
$ws.Range("B2:J29").Value = "UN"
$ws.Range("B30:J31").Value = "UN"
$ws.Range("B32:J33").Value = "UN"
